# Fruta / hortaliza, semanal
# Insert two new weekly records (date 2023-07-28 / serial 45135) for
# "Terminal La Palmera de La Serena" - Chirimoya, "Provincia de Limarí",
# right above the existing block that starts at row 214. This pushes the
# former rows 214-235 down to 216-237 (matching the target dimension
# A1:T237) and fills the two freshly inserted rows with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 214 (each Insert() shifts everything at/below
# row 214 down by one row, carrying formatting from the row above with it -
# which is exactly why the new D214/D215 already pick up the date style).
$ws.Rows.Item(214).Insert()
$ws.Rows.Item(214).Insert()

# New row 214: Primera
$ws.Range("A214").Value = 8
$ws.Range("B214").Value = "Terminal La Palmera de La Serena"
$ws.Range("C214").Value = "Coquimbo"
$ws.Range("D214").Value = 45135
$ws.Range("E214").Value = 4
$ws.Range("F214").Value = "Fruta"
$ws.Range("G214").Value = 100107
$ws.Range("H214").Value = "Otros"
$ws.Range("I214").Value = 100107002
$ws.Range("J214").Value = "Chirimoya"
$ws.Range("K214").Value = "Cultivar IV Región"
$ws.Range("L214").Value = "Primera"
$ws.Range("M214").Value = 500
$ws.Range("N214").Value = 25000
$ws.Range("O214").Value = 26000
$ws.Range("P214").Value = 25500
$ws.Range("Q214").Value = "$/bandeja 10 kilos"
$ws.Range("R214").Value = "Provincia de Limarí"
$ws.Range("S214").Value = 2550
$ws.Range("T214").Value = 10

# New row 215: Segunda
$ws.Range("A215").Value = 8
$ws.Range("B215").Value = "Terminal La Palmera de La Serena"
$ws.Range("C215").Value = "Coquimbo"
$ws.Range("D215").Value = 45135
$ws.Range("E215").Value = 4
$ws.Range("F215").Value = "Fruta"
$ws.Range("G215").Value = 100107
$ws.Range("H215").Value = "Otros"
$ws.Range("I215").Value = 100107002
$ws.Range("J215").Value = "Chirimoya"
$ws.Range("K215").Value = "Cultivar IV Región"
$ws.Range("L215").Value = "Segunda"
$ws.Range("M215").Value = 320
$ws.Range("N215").Value = 22000
$ws.Range("O215").Value = 23000
$ws.Range("P215").Value = 22500
$ws.Range("Q215").Value = "$/bandeja 10 kilos"
$ws.Range("R215").Value = "Provincia de Limarí"
$ws.Range("S215").Value = 2250
$ws.Range("T215").Value = 10
